$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("L1").Value = "break_on_off"

# Default every data row (2-73) in column L (break_on_off) to 0
for ($row = 2; $row -le 73; $row++) {
    $ws.Cells.Item($row, 12).Value = 0
}

# Rows where a break screen occurred (break_on_off = 1)
$breakRows = @(19, 37, 54)
foreach ($row in $breakRows) {
    $ws.Cells.Item($row, 12).Value = 1
}

[void]$ws.Range("L1:L73").Select()
